$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted as a new data row right above
# the existing row 91, pushing every subsequent row down by one.
$ws.Rows.Item(91).Insert()

$newRow = 91

$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = 44512
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100102
$ws.Cells.Item($newRow, 8).Value = "Cítricos"
$ws.Cells.Item($newRow, 9).Value = 100102006
$ws.Cells.Item($newRow, 10).Value = "Pomelo"
$ws.Cells.Item($newRow, 11).Value = "Start Ruby"
$ws.Cells.Item($newRow, 12).Value = "Primera"
$ws.Cells.Item($newRow, 13).Value = 200
$ws.Cells.Item($newRow, 14).Value = 11000
$ws.Cells.Item($newRow, 15).Value = 12000
$ws.Cells.Item($newRow, 16).Value = 11500
$ws.Cells.Item($newRow, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item($newRow, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($newRow, 19).Value = 821
$ws.Cells.Item($newRow, 20).Value = 14
